# Apply the edits described by the commit:
#   1. Slide 6's table switches from the custom "Table_0" style to the
#      built-in PowerPoint table style {6EC56CAB-2668-491E-816C-6C027E4FFA0D}.
#   2. The presentation's theme (ppt/theme/theme1.xml, used by the slide
#      master / all slides) is swapped from the "Integral" color scheme to
#      the stock "Office Theme" color scheme.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 --------------------------------------------
$s = $p.Slides.Item(6)
$tableShape = $s.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{6EC56CAB-2668-491E-816C-6C027E4FFA0D}")

# --- 2. Swap the theme colour scheme to the "Office Theme" palette -------
function ToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$scheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Length; $i++) {
    $scheme.Colors($i).RGB = ToOleColor($officeThemeColors[$i - 1])
}
